$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 1660
$ws.Range("I8").Value = 1660
$ws.Range("K8").Value = 4980
$ws.Range("M8").Value = -4841
# Row 15
$ws.Range("H15").Value = 265.34
$ws.Range("I15").Value = 265.34
$ws.Range("K15").Value = 796.02
$ws.Range("M15").Value = -627.02
# Row 111
$ws.Range("H111").Value = 3638.1667
$ws.Range("I111").Value = 5264.5
$ws.Range("J111").Value = 2825
$ws.Range("K111").Value = 15793.5
$ws.Range("L111").Value = 8475
$ws.Range("M111").Value = -12726.5
$ws.Range("N111").Value = -14609
# Row 132
$ws.Range("H132").Value = 3955.8518
$ws.Range("I132").Value = 3926.3044
$ws.Range("J132").Value = 4125.75
$ws.Range("K132").Value = 11778.9132
$ws.Range("L132").Value = 12377.25
$ws.Range("M132").Value = -9248.913199999999
$ws.Range("N132").Value = -17437.25
# Row 138
$ws.Range("H138").Value = 3714.3635
$ws.Range("I138").Value = 3118.913
$ws.Range("J138").Value = 3967.9814
$ws.Range("K138").Value = 9356.739
$ws.Range("L138").Value = 11903.9442
$ws.Range("M138").Value = -4216.739
$ws.Range("N138").Value = -22183.9442

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6108051.5
$ws.Range("I32").Value = 6858277
$ws.Range("J32").Value = 22888.889
$ws.Range("K32").Value = 6858277
$ws.Range("L32").Value = 22888.889
$ws.Range("M32").Value = -6857990
$ws.Range("N32").Value = -23462.889
# Row 97
$ws.Range("H97").Value = 849.8387
$ws.Range("I97").Value = 602.0454999999999
$ws.Range("J97").Value = 1455.5555
$ws.Range("K97").Value = 602.0454999999999
$ws.Range("L97").Value = 1455.5555
$ws.Range("M97").Value = -106.0454999999999
$ws.Range("N97").Value = -2447.5555
# Row 122
$ws.Range("H122").Value = 168425
$ws.Range("I122").Value = 251250
$ws.Range("J122").Value = 2775
$ws.Range("K122").Value = 753750
$ws.Range("L122").Value = 8325
$ws.Range("M122").Value = -751300
$ws.Range("N122").Value = -13225
# Row 132
$ws.Range("H132").Value = 1186262.4
$ws.Range("I132").Value = 2189.976
$ws.Range("J132").Value = 3348481.5
$ws.Range("K132").Value = 6569.928
$ws.Range("L132").Value = 10045444.5
$ws.Range("M132").Value = -4039.928
$ws.Range("N132").Value = -10050504.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 108
$ws.Range("H108").Value = 28333.334
$ws.Range("J108").Value = 28333.334
$ws.Range("L108").Value = 28333.334
$ws.Range("N108").Value = -36013.334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 428.14285
$ws.Range("I22").Value = 295.16666
$ws.Range("J22").Value = 527.875
$ws.Range("K22").Value = 295.16666
$ws.Range("L22").Value = 527.875
$ws.Range("M22").Value = 54.83334000000002
$ws.Range("N22").Value = -1227.875
# Row 31
$ws.Range("H31").Value = 5152.6045
$ws.Range("I31").Value = 1497.8064
$ws.Range("J31").Value = 7040.9165
$ws.Range("K31").Value = 1497.8064
$ws.Range("L31").Value = 7040.9165
$ws.Range("M31").Value = -1202.8064
$ws.Range("N31").Value = -7630.9165
# Row 34
$ws.Range("H34").Value = 5152.6045
$ws.Range("I34").Value = 1497.8064
$ws.Range("J34").Value = 7040.9165
$ws.Range("K34").Value = 1497.8064
$ws.Range("L34").Value = 7040.9165
$ws.Range("M34").Value = -1295.8064
$ws.Range("N34").Value = -7444.9165
# Row 132
$ws.Range("H132").Value = 16261916
$ws.Range("I132").Value = 16130533
$ws.Range("J132").Value = 16669203
$ws.Range("K132").Value = 48391599
$ws.Range("L132").Value = 50007609
$ws.Range("M132").Value = -48389069
$ws.Range("N132").Value = -50012669
# Row 134
$ws.Range("H134").Value = 7818182.5
$ws.Range("I134").Value = 11911386
$ws.Range("J134").Value = 3884.818
$ws.Range("K134").Value = 35734158
$ws.Range("L134").Value = 11654.454
$ws.Range("M134").Value = -35731623
$ws.Range("N134").Value = -16724.454
# Row 141
$ws.Range("H141").Value = 101707.22
$ws.Range("J141").Value = 97926.734
$ws.Range("L141").Value = 97926.734
$ws.Range("N141").Value = -108286.734

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1134
$ws.Range("I5").Value = 818.5925999999999
$ws.Range("J5").Value = 1582.2106
$ws.Range("K5").Value = 2455.7778
$ws.Range("L5").Value = 4746.6318
$ws.Range("M5").Value = -2343.7778
$ws.Range("N5").Value = -4970.6318
# Row 55
$ws.Range("H55").Value = 1764.6428
$ws.Range("J55").Value = 1792.6923
$ws.Range("L55").Value = 5378.0769
$ws.Range("N55").Value = -5732.0769
# Row 113
$ws.Range("H113").Value = 728.95
$ws.Range("I113").Value = 795
$ws.Range("J113").Value = 674.9091
$ws.Range("K113").Value = 2385
$ws.Range("L113").Value = 2024.7273
$ws.Range("M113").Value = -215
$ws.Range("N113").Value = -6364.7273
# Row 135
$ws.Range("H135").Value = 1134
$ws.Range("I135").Value = 818.5925999999999
$ws.Range("J135").Value = 1582.2106
$ws.Range("K135").Value = 7367.3334
$ws.Range("L135").Value = 14239.8954
$ws.Range("M135").Value = -4832.3334
$ws.Range("N135").Value = -19309.8954
# Row 140
$ws.Range("H140").Value = 1602.9
$ws.Range("I140").Value = 1151.2
$ws.Range("J140").Value = 2958
$ws.Range("K140").Value = 3453.6
$ws.Range("L140").Value = 8874
$ws.Range("M140").Value = 1726.4
$ws.Range("N140").Value = -19234

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 15000143
$ws.Range("I14").Value = 15000143
$ws.Range("K14").Value = 15000143
$ws.Range("M14").Value = -14999975
# Row 97
$ws.Range("H97").Value = 3241.2
$ws.Range("I97").Value = 3086
$ws.Range("K97").Value = 3086
$ws.Range("M97").Value = -2590
# Row 102
$ws.Range("H102").Value = 1814.5518
$ws.Range("I102").Value = 1812.32
$ws.Range("J102").Value = 1828.5
$ws.Range("K102").Value = 1812.32
$ws.Range("L102").Value = 1828.5
$ws.Range("M102").Value = -190.3199999999999
$ws.Range("N102").Value = -5072.5
# Row 113
$ws.Range("H113").Value = 93653.836
$ws.Range("I113").Value = 111913.3
$ws.Range("K113").Value = 111913.3
$ws.Range("M113").Value = -109743.3
# Row 122
$ws.Range("H122").Value = 3150
$ws.Range("I122").Value = 2300
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 6900
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -4450
$ws.Range("N122").Value = -16900

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2831.8572
$ws.Range("I132").Value = 1877.909
$ws.Range("K132").Value = 5633.727000000001
$ws.Range("M132").Value = -3103.727000000001
